$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 5143
$ws1.Range("F5").Value = 5143
$ws1.Range("F11").Value = 175
$ws1.Range("F12").Value = 8539
$ws1.Range("F13").Value = 8539
$ws1.Range("F17").Value = 2562
$ws1.Range("F19").Value = 2319
$ws1.Range("F25").Value = 6469
$ws1.Range("F27").Value = 73
$ws1.Range("F31").Value = 6980
$ws1.Range("F33").Value = 35
$ws1.Range("F43").Value = 2539
$ws1.Range("F46").Value = 1132
$ws1.Range("F47").Value = 66
$ws1.Range("F48").Value = 534
$ws1.Range("F49").Value = 2787
$ws1.Range("F50").Value = 86

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 189
$ws2.Range("F6").Value = 78
$ws2.Range("F15").Value = 28

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 5143
$ws4.Range("F4").Value = 5143
$ws4.Range("F10").Value = 175
$ws4.Range("F11").Value = 8539
$ws4.Range("F12").Value = 8539
$ws4.Range("F16").Value = 2562
$ws4.Range("F18").Value = 189
$ws4.Range("F20").Value = 2319
$ws4.Range("F21").Value = 78
$ws4.Range("F26").Value = 6469
$ws4.Range("F29").Value = 73
$ws4.Range("F33").Value = 6980
$ws4.Range("F34").Value = 35
$ws4.Range("F37").Value = 109
$ws4.Range("F41").Value = 2539
$ws4.Range("F43").Value = 1132
$ws4.Range("F44").Value = 66
$ws4.Range("F45").Value = 534
$ws4.Range("F47").Value = 2787
$ws4.Range("F48").Value = 86
$ws4.Range("F51").Value = 28
